$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 856, pushing existing rows 856..970 down to 857..971.
$ws.Rows("856:856").Insert()

# Populate the newly inserted row 856 with the new record's values.
$ws.Range("A856").Value = 8
$ws.Range("B856").Value = "Terminal La Palmera de La Serena"
$ws.Range("C856").Value = "Coquimbo"
$ws.Range("D856").Value = 45154
$ws.Range("E856").Value = 4
$ws.Range("F856").Value = 100112045
$ws.Range("G856").Value = "Zapallo"
$ws.Range("H856").Value = "Camote"
$ws.Range("I856").Value = "1a (guarda)"
$ws.Range("J856").Value = 1600
$ws.Range("K856").Value = 800
$ws.Range("L856").Value = 900
$ws.Range("M856").Value = 850
$ws.Range("N856").Value = '$/kilo (volumen en unidades)'
$ws.Range("O856").Value = "Región de O'Higgins"
$ws.Range("P856").Value = 850
$ws.Range("Q856").Value = 1
$ws.Range("R856").Value = "Hortaliza"
